$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the match name in C2 (it becomes a single blank space)
$ws.Range("C2").Value = " "

# 2. Add two new columns: PPDA 1st Half (DD) and PPDA 2nd Half (DE)
$ws.Range("DD1").Value = "PPDA 1st Half"
$ws.Range("DE1").Value = "PPDA 2nd Half"

$ddValues = @(5.95, 11.9, 5.64, 7.94, 4.24, 13.5, 8.64, 19.75, 6.79, 7.13, 5.9, 11.4, 14.4, 12.5, 14.54, 6.94, 5.85, 8.39, 5.27, 35.17, 6.36, 11.07, 8.47, 14.17, 9.89, 4.54, 13.8, 10.77, 18.5, 11.17)
$deValues = @(11.56, 10, 3.77, 6.83, 8.55, 5.27, 9.77, 8.6, 9.9, 3.95, 6.19, 27.17, 7.76, 13.91, 7.32, 4.82, 4.24, 12.08, 8.5, 9.07, 7.58, 8.15, 4.62, 6, 9.29, 7.29, 12, 10.27, 3.68, 7.21)

for ($i = 0; $i -lt 30; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 108).Value = $ddValues[$i]
    $ws.Cells.Item($row, 109).Value = $deValues[$i]
}

# 3. Column widths for the new columns (best-fit approximation)
$ws.Range("DD1").ColumnWidth = 11.7
$ws.Range("DE1").ColumnWidth = 12.25

# 4. Restore the selection left by the editor
[void]$ws.Range("CJ1").Select()
[void]$ws.Range("DH16").Select()
